$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; this shifts the existing rows 5-84 down
# to become rows 6-85 (matching the dimension growing from A1:T84 to A1:T85).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value  = 5
$ws.Cells.Item(5, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value  = "Maule"
$ws.Cells.Item(5, 4).Value  = "2021-10-27"
$ws.Cells.Item(5, 5).Value  = 7
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100108
$ws.Cells.Item(5, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(5, 9).Value  = 100108002
$ws.Cells.Item(5, 10).Value = "Mango"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 210
$ws.Cells.Item(5, 14).Value = 8000
$ws.Cells.Item(5, 15).Value = 8000
$ws.Cells.Item(5, 16).Value = 8000
$ws.Cells.Item(5, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(5, 18).Value = "Perú"
$ws.Cells.Item(5, 19).Value = 2000
$ws.Cells.Item(5, 20).Value = 4
